$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 79, shifting existing rows 79-140 down to 84-145.
$ws.Range("A79:B83").EntireRow.Insert()

# Populate the newly inserted rows with the new vocabulary/phrase entries.
$ws.Range("A79").Value = "Thank you for everything."
$ws.Range("B79").Value = "いろいろおせわになりました。"

$ws.Range("A80").Value = "Please take care of yourself."
$ws.Range("B80").Value = "体に気をつけてください。|からだにきをつけてください。"

$ws.Range("A81").Value = "I am looking forward to seeing you."
$ws.Range("B81").Value = "お会いできるのを楽しみにしています。|おあいできるをたのしみにしています。"

$ws.Range("A82").Value = "Congratulations on..."
$ws.Range("B82").Value = "～おめでとう（ございます）。"

$ws.Range("A83").Value = "Happy Birthday."
$ws.Range("B83").Value = "（お）たんじょうびおめでとう。"
